{"js": "// Remove the \"-Randomize level generation\" line (whole paragraph, including\n// its paragraph mark) and relocate the trailing \"_GoBack\" bookmark so it\n// sits at the start of the \"-Add to list of possible levels\" paragraph\n// instead of its own now-empty trailing paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Step 1: delete the \"-Randomize level generation\" paragraph entirely ---\nconst toRemove = paragraphs.items.find(\n  (p) => p.text === \"-Randomize level generation\"\n);\nif (toRemove) {\n  toRemove.delete();\n  await context.sync();\n}\n\n// --- Step 2: locate the \"-Add to list of possible levels\" paragraph ---\nconst remaining = body.paragraphs;\nremaining.load(\"items/text\");\nawait context.sync();\n\nconst addParagraph = remaining.items.find(\n  (p) => p.text === \"-Add to list of possible levels\"\n);\n\n// --- Step 3: move the _GoBack bookmark to the start of that paragraph ---\ncontext.document.deleteBookmark(\"_GoBack\");\nif (addParagraph) {\n  const startRange = addParagraph.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "# Remove the \"-Randomize level generation\" line (whole paragraph, including\n# its paragraph mark) and relocate the trailing \"_GoBack\" bookmark so it sits\n# at the start of the \"-Add to list of possible levels\" paragraph instead of\n# its own now-empty trailing paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"-Randomize level generation\" paragraph entirely ---\n$delRng = $d.Content\n$null = $delRng.Find.Execute('-Randomize level generation')\n$delPara = $delRng.Paragraphs(1).Range\n$delPara.Delete()\n\n# --- Step 2: locate the \"-Add to list of possible levels\" paragraph ---\n$addRng = $d.Content\n$null = $addRng.Find.Execute('-Add to list of possible levels')\n$addPara = $addRng.Paragraphs(1).Range\n\n# --- Step 3: move the _GoBack bookmark to the start of that paragraph ---\nif ($d.Bookmarks.Exists('_GoBack')) {\n    $d.Bookmarks('_GoBack').Delete()\n}\n\n$target = $addPara.Duplicate\n$target.Collapse(1)\n$d.Bookmarks.Add('_GoBack', $target)\n"}
